$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.748.90'
$ws.Range('E2').Value = '  +1.78%  '
$ws.Range('D3').Value = '3.562.30'
$ws.Range('E3').Value = '  +1.63%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '583.42'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.47%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '189.54'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.69%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.624'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +2.02%  '
$ws.Range('D8').Value = '3.553.82'
$ws.Range('E8').Value = '  +1.58%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('E10').Value = '  +15.77%  '
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '54.60'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('E13').Value = '  +6.02%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.51'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.72%  '
$ws.Range('D15').Value = '4.131.35'
$ws.Range('E15').Value = '  +1.63%  '
$ws.Range('D16').Value = '70.742.12'
$ws.Range('E16').Value = '  +2.01%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.572.61'
$ws.Range('E17').Value = '  +1.85%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '19.15'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.77'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +4.24%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '569.56'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +5.40%  '
$ws.Range('E21').Value = '  +0.78%  '
$ws.Range('E22').Value = '  -0.70%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '17.93'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.59%  '
$ws.Range('E24').Value = '  +4.15%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '4.90'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.33%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '94.19'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.60%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.15'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.17%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.93'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.51%  '
$ws.Range('E29').Value = '  +1.86%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '32.53'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.27%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.19'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.91%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '12.30'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.67%  '
$ws.Range('E33').Value = '  +2.18%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '63.85'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.13%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.73'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +22.19%  '
$ws.Range('E36').Value = '  +5.25%  '
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '532.52'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.56%  '
$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.411'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.50%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '38.40'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.12%  '
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '3.630.93'
$ws.Range('E41').Value = '  +10.25%  '
$ws.Range('D42').Value = '0.0₃0798'
$ws.Range('E42').Value = '  +4.41%  '
$ws.Range('E43').Value = '  +5.10%  '
$ws.Range('E44').Value = '  +1.56%  '
$ws.Range('E45').Value = '  +5.22%  '
$ws.Range('E46').Value = '  -1.17%  '
$ws.Range('E47').Value = '  -0.92%  '
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.28'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +4.19%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.138'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.24%  '
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('E51').Value = '  +6.80%  '
